$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, pushing existing rows 149:212 down to 150:213.
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new record's data.
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44917
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = 100112030
$ws.Cells.Item(149, 7).Value = "Poroto granado"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 73
$ws.Cells.Item(149, 11).Value = 28000
$ws.Cells.Item(149, 12).Value = 30000
$ws.Cells.Item(149, 13).Value = 28959
$ws.Cells.Item(149, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(149, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(149, 16).Value = 1158
$ws.Cells.Item(149, 17).Value = 25
$ws.Cells.Item(149, 18).Value = "Hortaliza"
